$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Brochage v2")

# --- Row 7 / Row 8 in column A: move into the "orange" colour block ---
# (matches cells A11/A12 which already carry that fill) by copying their
# format onto A7/A8 so the existing style entry is reused rather than a
# brand new one being minted.
$ws2.Range("A11").Copy() | Out-Null
$ws2.Range("A7:A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Rotate the CE / CODD_A / CODD_B labels (and their highlight) down one row ---
# Before: A22=CODD_B (plain), A23=CODD_A (plain), A24=CE (blue highlight)
# After : A22=CE (blue highlight), A23=CODD_B (plain), A24=CODD_A (plain)
$ws2.Range("A13").Copy() | Out-Null
$ws2.Range("A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws2.Range("A6").Copy() | Out-Null
$ws2.Range("A24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws2.Range("A22").Value = "CE"
$ws2.Range("A23").Value = "CODD_B"
$ws2.Range("A24").Value = "CODD_A"

# --- Highlight E7 (empty note cell) in bold red ---
$ws2.Range("E7").Font.Bold = $true

# --- Leave the selection where the author last clicked ---
$ws2.Activate()
$ws2.Range("B32").Select()
